# Ajuste do cadastro de pessoas em massa, para enviar convite por e-mail.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove existing hyperlinks (e-mail cells become plain values)
$ws.Hyperlinks.Delete()

# Row 2: Toni Garrido
$ws.Range("A2").Value = "Toni Garrido"
$ws.Range("B2").Value = "técnico"
$ws.Range("C2").Value = "renato@ispn.org.br"
$ws.Range("D2").Value = 666666
$ws.Range("E2").Value = "cod/04"

# Row 3: Toni Beloto
$ws.Range("A3").Value = "Toni Beloto"
$ws.Range("B3").Value = "técnico"
$ws.Range("C3").Value = "jenipapos@yahoo.com.br"
$ws.Range("D3").Value = 55555

# Row 4 (new): Paulo Miklos
$ws.Range("A4").Value = "Paulo Miklos"
$ws.Range("B4").Value = "técnico"
$ws.Range("C4").Value = "renato@renato.org.br"
$ws.Range("D4").Value = 44444

# A few blank formatted rows below the data, extending the used range to row 7
$ws.Range("A5").Borders.LineStyle = 0
$ws.Range("A6").Borders.LineStyle = 0
$ws.Range("A7").Borders.LineStyle = 0

# Move selection to B13, matching the saved cursor position in the source file
$ws.Range("B13").Select()
